{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the paragraph that ends with \"datepicker\" \u2014 the new content is\n// inserted right after it.\nlet datepickerPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"datepicker\") !== -1) {\n    datepickerPara = p;\n  }\n}\n\nif (datepickerPara === null) {\n  throw new Error(\"Could not locate the 'datepicker' paragraph\");\n}\n\nconst noteText =\n  \"\u043a\u0430\u0442\u043e \u0441\u0435 \u0440\u0435\u0433\u0438\u0441\u0442\u0440\u0438\u0440\u0430\u043c\u0435 \u0442\u0440 \u0438\u043b\u0438 \u0434\u0430 \u0437\u0430\u0440\u0435\u0436\u0434\u0430\u043c\u0435 \u0434\u0430\u0448\u0431\u043e\u0440\u0434\u0430 \u0438\u043b\u0438 \u0434\u0430 \u0441\u0435 \u043e\u0442\u043f\u0438\u0441\u0432\u0430 \" +\n  \"\u043f\u043e\u0442\u0440\u0435\u0431\u0438\u0442\u0435\u043b\u044f \u0438 \u0434\u0430 \u043f\u0440\u0430\u0449\u0430 \u043d\u0430 \u0433\u043d\u0430\u0447\u0430\u043b\u043d\u0430\u0442\u0430 \u0441\u0442\u0440\u0430\u043d\u0438\u0446\u0430\";\n\n// Build the four new paragraphs (two blank, one with the note, one blank\n// that will carry the relocated \"_GoBack\" bookmark) as a flat-OPC OOXML\n// fragment so they land as genuinely empty <w:p/> elements instead of\n// picking up a stray empty run.\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p/>\" +\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>\" + noteText + \"</w:t></w:r></w:p>\" +\n  '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst afterDatepicker = datepickerPara.getRange(\"After\");\nafterDatepicker.insertOoxml(flatOpc, \"After\");\nawait context.sync();\n\n// The \"_GoBack\" bookmark used to sit at the end of the \"\u0417\u0430\u0449\u0438\u0442\u0430 \u043d\u0430 \u0411\u0414\"\n// paragraph; now that it has been re-created on the new trailing empty\n// paragraph, drop the original one so the name stays unique.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that ends with \"datepicker\" \u2014 the new content is\n# inserted right after it.\n$datepickerPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -match \"datepicker\") {\n        $datepickerPara = $para\n    }\n}\n\nif ($datepickerPara -eq $null) {\n    throw \"Could not locate the 'datepicker' paragraph\"\n}\n\n$noteText = \"\u043a\u0430\u0442\u043e \u0441\u0435 \u0440\u0435\u0433\u0438\u0441\u0442\u0440\u0438\u0440\u0430\u043c\u0435 \u0442\u0440 \u0438\u043b\u0438 \u0434\u0430 \u0437\u0430\u0440\u0435\u0436\u0434\u0430\u043c\u0435 \u0434\u0430\u0448\u0431\u043e\u0440\u0434\u0430 \u0438\u043b\u0438 \u0434\u0430 \u0441\u0435 \u043e\u0442\u043f\u0438\u0441\u0432\u0430 \u043f\u043e\u0442\u0440\u0435\u0431\u0438\u0442\u0435\u043b\u044f \u0438 \u0434\u0430 \u043f\u0440\u0430\u0449\u0430 \u043d\u0430 \u0433\u043d\u0430\u0447\u0430\u043b\u043d\u0430\u0442\u0430 \u0441\u0442\u0440\u0430\u043d\u0438\u0446\u0430\"\n\n# Insert the four new paragraphs (two blank, one with the note, one blank\n# that will carry the relocated \"_GoBack\" bookmark) as a single flat-OPC\n# OOXML fragment at a collapsed range right after the paragraph mark of\n# \"datepicker\" \u2014 this yields genuine empty <w:p/> elements instead of\n# leaving a stray empty run behind.\n$insertAt = $d.Range($datepickerPara.Range.End, $datepickerPara.Range.End)\n$xml = '<?xml version=\"1.0\" standalone=\"yes\"?>' `\n    + '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' `\n    + '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' `\n    + '<pkg:xmlData>' `\n    + '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' `\n    + '<w:body>' `\n    + '<w:p/>' `\n    + '<w:p/>' `\n    + '<w:p><w:r><w:t>' + $noteText + '</w:t></w:r></w:p>' `\n    + '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' `\n    + '</w:body>' `\n    + '</w:document>' `\n    + '</pkg:xmlData>' `\n    + '</pkg:part>' `\n    + '</pkg:package>'\n$insertAt.InsertXML($xml)\n\n# The \"_GoBack\" bookmark used to sit at the end of the \"\u0417\u0430\u0449\u0438\u0442\u0430 \u043d\u0430 \u0411\u0414\"\n# paragraph; now that it has been re-created on the new trailing empty\n# paragraph, drop the original one so the name stays unique.\n$d.Bookmarks(\"_GoBack\").Delete()\n"}
